$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1104
$ws.Range("I2").Value = 2798
$ws.Range("J2").Value = 11437
$ws.Range("K2").Value = 67
$ws.Range("L2").Value = 3120
$ws.Range("M2").Value = 192
$ws.Range("N2").Value = 2019
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 35
$ws.Range("Q2").Value = 14
$ws.Range("R2").Value = 154
$ws.Range("S2").Value = 1200
$ws.Range("T2").Value = 1941
$ws.Range("U2").Value = 135
$ws.Range("V2").Value = 17677
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 17738
$ws.Range("Y2").Value = 36
$ws.Range("Z2").Value = 255
$ws.Range("AA2").Value = 111
